# Add a new "E-mail do Comprador" column right after the existing
# "Comprador" column (column K) in the order template worksheet.
#
# This inserts a new column at position L (12), shifting every
# subsequent column one place to the right, and fills in the header
# text for the newly inserted column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column L is the 12th column - insert a blank column there, pushing
# the former column L ("CPF/CNPJ") and everything after it right by one.
$ws.Columns.Item(12).Insert()

# Set the header text for the newly inserted column.
$ws.Cells.Item(1, 12).Value2 = "E-mail do Comprador"
